# Update "想去人数" (number of people interested) values on the
# "展览" and "全部类型" sheets to reflect the newly scraped counts.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 11590
    $ws.Range("F3").Value = 11124
    $ws.Range("F17").Value = 43
    $ws.Range("F21").Value = 10880
}
